$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1229.1
$ws.Range("J9").Value = 1815.3334
$ws.Range("L9").Value = 1815.3334
$ws.Range("N9").Value = -2153.3334
$ws.Range("H12").Value = 1979.7
$ws.Range("J12").Value = 2642.4285
$ws.Range("L12").Value = 2642.4285
$ws.Range("N12").Value = -2982.4285
$ws.Range("H33").Value = 652.9474
$ws.Range("J33").Value = 4499.5
$ws.Range("L33").Value = 4499.5
$ws.Range("N33").Value = -4957.5
$ws.Range("H49").Value = 6510
$ws.Range("J49").Value = 6510
$ws.Range("L49").Value = 19530
$ws.Range("N49").Value = -19802
$ws.Range("H98").Value = 578.5
$ws.Range("I98").Value = 556.3158
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 556.3158
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 941.6842
$ws.Range("N98").Value = -3996
$ws.Range("H122").Value = 578.5
$ws.Range("I122").Value = 556.3158
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 1668.9474
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 781.0526
$ws.Range("N122").Value = -7900
$ws.Range("H138").Value = 3497.4893
$ws.Range("I138").Value = 3024.875
$ws.Range("J138").Value = 3741.4194
$ws.Range("K138").Value = 9074.625
$ws.Range("L138").Value = 11224.2582
$ws.Range("M138").Value = -3934.625
$ws.Range("N138").Value = -21504.2582

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 2016.3334
$ws.Range("I54").Value = 2016.3334
$ws.Range("K54").Value = 2016.3334
$ws.Range("M54").Value = -1532.3334
$ws.Range("H99").Value = 1415.4
$ws.Range("I99").Value = 1363.9231
$ws.Range("K99").Value = 1363.9231
$ws.Range("M99").Value = 134.0769
$ws.Range("H103").Value = 37223.43
$ws.Range("J103").Value = 37223.43
$ws.Range("L103").Value = 37223.43
$ws.Range("N103").Value = -39567.43
$ws.Range("H134").Value = 4080.5833
$ws.Range("I134").Value = 1834.8572
$ws.Range("K134").Value = 5504.571599999999
$ws.Range("M134").Value = -2969.571599999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 136.625
$ws.Range("I7").Value = 75.92308
$ws.Range("K7").Value = 75.92308
$ws.Range("M7").Value = 37.07692

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 467.75
$ws.Range("I46").Value = 457.5
$ws.Range("J46").Value = 498.5
$ws.Range("K46").Value = 1372.5
$ws.Range("L46").Value = 1495.5
$ws.Range("M46").Value = -1281.5
$ws.Range("N46").Value = -1677.5
$ws.Range("H86").Value = 1191.0769
$ws.Range("I86").Value = 205.8
$ws.Range("J86").Value = 1806.875
$ws.Range("K86").Value = 617.4000000000001
$ws.Range("L86").Value = 5420.625
$ws.Range("M86").Value = 568.5999999999999
$ws.Range("N86").Value = -7792.625
$ws.Range("H88").Value = 17008
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()   # was -44572
$ws.Range("H89").Value = 1191.0769
$ws.Range("I89").Value = 205.8
$ws.Range("J89").Value = 1806.875
$ws.Range("K89").Value = 1852.2
$ws.Range("L89").Value = 16261.875
$ws.Range("M89").Value = 4075.8
$ws.Range("N89").Value = -28117.875
$ws.Range("H91").Value = 17008
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()   # was -43518
$ws.Range("H132").Value = 6121.643
$ws.Range("I132").Value = 5386.143
$ws.Range("K132").Value = 48475.287
$ws.Range("M132").Value = -45945.287

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 745.2857
$ws.Range("I2").Value = 79.25
$ws.Range("K2").Value = 79.25
$ws.Range("M2").Value = 33.75
$ws.Range("H3").Value = 3156.4285
$ws.Range("J3").Value = 5750
$ws.Range("L3").Value = 5750
$ws.Range("N3").Value = -5982
$ws.Range("H10").Value = 9281.6
$ws.Range("J10").Value = 2503
$ws.Range("L10").Value = 2503
$ws.Range("N10").Value = -2841
$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 15000
$ws.Range("J43").Value = 25000
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 25000
$ws.Range("M43").Value = -14849
$ws.Range("N43").Value = -25302
$ws.Range("H55").Value = 765000
$ws.Range("I55").Value = 30000
$ws.Range("K55").Value = 30000
$ws.Range("M55").Value = -29673
$ws.Range("H64").Value = 79000
$ws.Range("J64").Value = 79000
$ws.Range("L64").Value = 79000
$ws.Range("N64").Value = -79496
$ws.Range("H67").Value = 79000
$ws.Range("J67").Value = 79000
$ws.Range("L67").Value = 79000
$ws.Range("N67").Value = -80716
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()   # was -46872
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()   # was -144360
$ws.Range("H105").Value = 20671
$ws.Range("J105").Value = 20671
$ws.Range("L105").Value = 20671
$ws.Range("N105").Value = -27659

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5863.3076
$ws.Range("J22").Value = 7221.222
$ws.Range("L22").Value = 7221.222
$ws.Range("N22").Value = -7811.222
$ws.Range("H27").Value = 5863.3076
$ws.Range("J27").Value = 7221.222
$ws.Range("L27").Value = 7221.222
$ws.Range("N27").Value = -7435.222
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()   # was -39450
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()   # was -40560
$ws.Range("H93").Value = 9338.4
$ws.Range("I93").Value = 2277.8
$ws.Range("J93").Value = 16399
$ws.Range("K93").Value = 2277.8
$ws.Range("L93").Value = 16399
$ws.Range("M93").Value = -1029.8
$ws.Range("N93").Value = -18895
$ws.Range("H132").Value = 4787.7646
$ws.Range("I132").Value = 2710.5557
$ws.Range("K132").Value = 8131.6671
$ws.Range("M132").Value = -5601.6671

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 683.62964
$ws.Range("J113").Value = 466
$ws.Range("L113").Value = 1398
$ws.Range("N113").Value = -5738
$ws.Range("H122").Value = 2195.8235
$ws.Range("I122").Value = 1718
$ws.Range("K122").Value = 5154
$ws.Range("M122").Value = -2704
$ws.Range("H132").Value = 7597.625
$ws.Range("I132").Value = 8041
$ws.Range("J132").Value = 7331.6
$ws.Range("K132").Value = 24123
$ws.Range("L132").Value = 21994.8
$ws.Range("M132").Value = -21593
$ws.Range("N132").Value = -27054.8
